$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "23.415.86"
$ws.Range("E2").Value = "  +0.73%  "
Set-TextValue $ws.Range("D3") "1.639.50"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.10%  "
Set-TextValue $ws.Range("D6") "304.89"
$ws.Range("E6").Value = "  +0.51%  "
Set-TextValue $ws.Range("D7") "0.3735"
$ws.Range("E7").Value = "  -1.03%  "
Set-TextValue $ws.Range("D8") "52.16"
$ws.Range("E8").Value = "  +0.89%  "
Set-TextValue $ws.Range("D9") "0.3628"
$ws.Range("E9").Value = "  -0.19%  "
Set-TextValue $ws.Range("D10") "1.249"
$ws.Range("E10").Value = "  -1.93%  "
Set-TextValue $ws.Range("D11") "0.08121"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  +0.08%  "
Set-TextValue $ws.Range("D13") "22.86"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -0.14%  "
Set-TextValue $ws.Range("D15") "0.00001271"
$ws.Range("E15").Value = "  +1.71%  "
Set-TextValue $ws.Range("D16") "7.279"
$ws.Range("E16").Value = "  -1.78%  "
Set-TextValue $ws.Range("D17") "1.628.95"
$ws.Range("E17").Value = "  +1.46%  "
Set-TextValue $ws.Range("D18") "94.42"
$ws.Range("E18").Value = "  +0.53%  "
Set-TextValue $ws.Range("D19") "0.06877"
$ws.Range("E19").Value = "  -0.21%  "
Set-TextValue $ws.Range("D20") "18.14"
$ws.Range("E20").Value = "  -0.06%  "
Set-TextValue $ws.Range("D21") "6.516"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  +0.07%  "
Set-TextValue $ws.Range("D23") "23.415.58"
$ws.Range("E23").Value = "  +0.76%  "
Set-TextValue $ws.Range("D24") "12.75"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("E25").Value = "  +1.03%  "
Set-TextValue $ws.Range("D26") "3.022"
$ws.Range("E26").Value = "  +0.65%  "
Set-TextValue $ws.Range("D27") "21.21"
$ws.Range("E27").Value = "  -0.10%  "
Set-TextValue $ws.Range("D28") "151.95"
$ws.Range("E28").Value = "  +1.23%  "
Set-TextValue $ws.Range("D29") "5.287"
$ws.Range("E29").Value = "  +0.70%  "
Set-TextValue $ws.Range("D30") "134.90"
$ws.Range("E30").Value = "  +0.53%  "
Set-TextValue $ws.Range("D31") "2.293"
$ws.Range("E31").Value = "  -3.71%  "
Set-TextValue $ws.Range("D32") "1.808.78"
$ws.Range("E32").Value = "  +1.72%  "
Set-TextValue $ws.Range("D33") "6.770"
$ws.Range("E33").Value = "  +0.09%  "
Set-TextValue $ws.Range("D34") "0.9523"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("E35").Value = "  +3.52%  "
Set-TextValue $ws.Range("D36") "10.32"
$ws.Range("E36").Value = "  +0.57%  "
Set-TextValue $ws.Range("D37") "0.2525"
$ws.Range("E37").Value = "  -0.33%  "
Set-TextValue $ws.Range("D38") "0.07209"
$ws.Range("E38").Value = "  -4.31%  "
Set-TextValue $ws.Range("D39") "0.08784"
$ws.Range("E39").Value = "  -0.28%  "
Set-TextValue $ws.Range("D40") "6.068"
$ws.Range("E40").Value = "  -0.47%  "
Set-TextValue $ws.Range("D41") "1.376"
$ws.Range("E41").Value = "  -1.52%  "
Set-TextValue $ws.Range("D42") "0.7044"
$ws.Range("E42").Value = "  -1.03%  "
Set-TextValue $ws.Range("D43") "12.44"
$ws.Range("E43").Value = "  -0.73%  "
Set-TextValue $ws.Range("D44") "15.96"
$ws.Range("E44").Value = "  +2.14%  "
Set-TextValue $ws.Range("D45") "0.6499"
$ws.Range("E45").Value = "  -0.75%  "
Set-TextValue $ws.Range("D46") "2.329"
Set-TextValue $ws.Range("D47") "1.000"
$ws.Range("E47").Value = "  +0.10%  "
Set-TextValue $ws.Range("D48") "4.008"
$ws.Range("E48").Value = "  -0.17%  "
Set-TextValue $ws.Range("D49") "0.07965"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  -2.97%  "
Set-TextValue $ws.Range("D51") "1.197"
$ws.Range("E51").Value = "  -0.65%  "
